$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are stored as text so that numeric-looking
# strings (e.g. "1.002", "0.00001227") are preserved exactly as text and are
# not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.036.35"
$ws.Range("E2").Value = "  -0.63%  "

$ws.Range("D3").Value = "1.590.82"
$ws.Range("E3").Value = "  -0.64%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("D6").Value = "301.52"
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").Value = "0.3769"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "0.3597"
$ws.Range("E8").Value = "  -1.75%  "

$ws.Range("D9").Value = "50.47"
$ws.Range("E9").Value = "  +5.91%  "

$ws.Range("D10").Value = "1.003"
$ws.Range("E10").Value = "  +0.14%  "

$ws.Range("D11").Value = "1.227"
$ws.Range("E11").Value = "  -4.07%  "

$ws.Range("D12").Value = "0.08069"
$ws.Range("E12").Value = "  -0.22%  "

$ws.Range("D13").Value = "22.03"
$ws.Range("E13").Value = "  -4.15%  "

$ws.Range("D14").Value = "6.485"
$ws.Range("E14").Value = "  -2.32%  "

$ws.Range("D15").Value = "7.262"
$ws.Range("E15").Value = "  -4.25%  "

$ws.Range("D16").Value = "0.00001227"
$ws.Range("E16").Value = "  -3.17%  "

$ws.Range("D17").Value = "1.592.12"
$ws.Range("E17").Value = "  -0.28%  "

$ws.Range("D18").Value = "92.42"
$ws.Range("E18").Value = "  +0.88%  "

$ws.Range("D19").Value = "0.06831"

$ws.Range("D20").Value = "17.97"
$ws.Range("E20").Value = "  -2.65%  "

$ws.Range("D21").Value = "6.463"
$ws.Range("E21").Value = "  -2.12%  "

$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").Value = "12.90"
$ws.Range("E23").Value = "  -1.23%  "

$ws.Range("D24").Value = "23.041.05"
$ws.Range("E24").Value = "  -0.62%  "

$ws.Range("D25").Value = "2.378"
$ws.Range("E25").Value = "  +0.75%  "

$ws.Range("D26").Value = "2.800"
$ws.Range("E26").Value = "  -3.86%  "

$ws.Range("D27").Value = "20.91"
$ws.Range("E27").Value = "  -1.00%  "

$ws.Range("D28").Value = "148.36"
$ws.Range("E28").Value = "  -1.87%  "

$ws.Range("D29").Value = "5.221"
$ws.Range("E29").Value = "  -0.45%  "

$ws.Range("D30").Value = "133.09"
$ws.Range("E30").Value = "  +0.79%  "

$ws.Range("D31").Value = "2.367"
$ws.Range("E31").Value = "  -3.10%  "

$ws.Range("D32").Value = "6.567"
$ws.Range("E32").Value = "  -7.78%  "

$ws.Range("D33").Value = "1.772.24"
$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("D34").Value = "0.9475"
$ws.Range("E34").Value = "  -3.48%  "

$ws.Range("D35").Value = "0.07411"
$ws.Range("E35").Value = "  -4.17%  "

$ws.Range("D36").Value = "10.10"
$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("D37").Value = "0.02681"
$ws.Range("E37").Value = "  -3.75%  "

$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "0.08793"
$ws.Range("E38").Value = "  -1.03%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "6.067"
$ws.Range("E39").Value = "  -4.12%  "

$ws.Range("D40").Value = "0.2483"
$ws.Range("E40").Value = "  -2.64%  "

$ws.Range("D41").Value = "1.357"
$ws.Range("E41").Value = "  -3.38%  "

$ws.Range("D42").Value = "0.6928"
$ws.Range("E42").Value = "  -3.30%  "

$ws.Range("D43").Value = "12.16"
$ws.Range("E43").Value = "  -4.86%  "

$ws.Range("D44").Value = "15.02"
$ws.Range("E44").Value = "  -5.98%  "

$ws.Range("D45").Value = "0.6476"
$ws.Range("E45").Value = "  -2.67%  "

$ws.Range("D46").Value = "3.999"
$ws.Range("E46").Value = "  +0.81%  "

$ws.Range("D47").Value = "2.258"
$ws.Range("E47").Value = "  -2.43%  "

$ws.Range("D48").Value = "131.15"
$ws.Range("E48").Value = "  -0.35%  "

$ws.Range("D49").Value = "0.07900"
$ws.Range("E49").Value = "  -1.12%  "

$ws.Range("D50").Value = "1.202"
$ws.Range("E50").Value = "  +2.47%  "

$ws.Range("D51").Value = "1.212"
$ws.Range("E51").Value = "  +3.80%  "
